$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 11, shifting existing rows 11:21 down to 12:22
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new Albahaca price record
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44494
$ws.Range("D11").NumberFormat = $ws.Range("D12").NumberFormat
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112052
$ws.Range("G11").Value = "Albahaca"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 2400
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2450
$ws.Range("N11").Value = "$/paquete"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 2450
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
